$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167 — this shifts the existing rows 167..269
# down to 168..270 (matching the diff, which shows every row from 167 to
# 269 taking on the values previously held by the row above it, and a new
# row 270 appearing with the data that used to be in row 269).
$ws.Rows(167).Insert()

# Populate the newly inserted row 167 with the new weekly price record.
$ws.Cells(167, 1).Value2 = 11
$ws.Cells(167, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells(167, 3).Value2 = "Bíobío"
$ws.Cells(167, 4).Value2 = 44673
$ws.Cells(167, 5).Value2 = 8
$ws.Cells(167, 6).Value2 = 100112023
$ws.Cells(167, 7).Value2 = "Brócoli"
$ws.Cells(167, 8).Value2 = "Sin especificar"
$ws.Cells(167, 9).Value2 = "Primera"
$ws.Cells(167, 10).Value2 = 2500
$ws.Cells(167, 11).Value2 = 750
$ws.Cells(167, 12).Value2 = 800
$ws.Cells(167, 13).Value2 = 780
$ws.Cells(167, 14).Value2 = "$/unidad"
$ws.Cells(167, 15).Value2 = "Región Metropolitana"
$ws.Cells(167, 16).Value2 = 780
$ws.Cells(167, 17).Value2 = 1
$ws.Cells(167, 18).Value2 = "Hortaliza"
